$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.6742952501365937
$ws.Range("D2").Value = 0.753548703160011

$ws.Range("B3").Value = 0.7029975329123593
$ws.Range("C3").Value = -0.7485016639922196
$ws.Range("D3").Value = 0.7117477281713626

$ws.Range("B4").Value = -0.5806511269897915
$ws.Range("C4").Value = -0.6910453398985321
$ws.Range("D4").Value = 0.6621674588708382

$ws.Range("B5").Value = -0.5736507905492713
$ws.Range("C5").Value = -0.6187236023588955
$ws.Range("D5").Value = -0.6650488329138234

$ws.Range("B6").Value = 0.5875158140643648
$ws.Range("C6").Value = -0.6759856675442053
$ws.Range("D6").Value = -0.7221743093942519

$ws.Range("B7").Value = 0.7683279318924957
$ws.Range("C7").Value = -0.6178857552030708
$ws.Range("D7").Value = 0.7397933964685425

$ws.Range("B8").Value = 0.7627723112583483
$ws.Range("C8").Value = 0.6158911077634686
$ws.Range("D8").Value = 0.789432957893222

$ws.Range("B9").Value = -0.6659068812185581
$ws.Range("C9").Value = 0.6543146566019556
$ws.Range("D9").Value = -0.6570529102178287
